$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: insert a new paragraph "Branch Hello!" right after the existing
# "Hello!" paragraph (same Body style / indent / Calibri formatting).
# ---------------------------------------------------------------------------
$helloPara = $d.Paragraphs(1)
$helloPara.Range.InsertParagraphAfter()
$branchPara = $d.Paragraphs(2)
$branchPara.Range.Text = "Branch Hello!"

# ---------------------------------------------------------------------------
# Change 2: the paragraph that reads
#   "You then must click the Download For Windows (64Bit) button to
#    download the .exe file for Windows operating systems. ..."
# originally has the first part split across three runs
#   [You then must click the Download ][For][ Windows (64Bit) button to
#   download the .exe file]
# with <w:proofErr type="gramStart"/> / <w:proofErr type="gramEnd"/> markers
# wrapping the middle run. Target: merge those three runs (and drop the
# proofErr markers) into a single run, while leaving the two runs that
# follow ("  for Windows operating systems" and ". Once downloaded, ...")
# untouched as separate runs.
# ---------------------------------------------------------------------------

# Locate the paragraph by its distinctive text (index shifted by the
# paragraph inserted above).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text.StartsWith("You then must click the Download")) {
        $target = $cand
        break
    }
}

$targetRange = $target.Range
$mergedText = "You then must click the Download For Windows (64Bit) button to download the .exe file"
$mergeLen = $mergedText.Length

# Step A: force a genuine text change across exactly the first three runs so
# the host collapses them (plus whatever trailing runs share the same
# proofErr-free run block) into a single run. Using a distinct marker
# guarantees the text differs from the original (an identity "replace"
# would be treated as a no-op and would not restructure the runs).
$startA = $targetRange.Start
$endA = $targetRange.Start + $mergeLen
$rngA = $d.Range($startA, $endA)
$rngA.Text = $mergedText + "@@MARK@@"

# Step B: strip the temporary marker back out, restoring the exact target
# text, without touching anything beyond it.
$target2 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text.StartsWith("You then must click the Download")) {
        $target2 = $cand
        break
    }
}
$targetRange2 = $target2.Range
$fullText2 = $targetRange2.Text
$markIdx = $fullText2.IndexOf("@@MARK@@")
$startMark = $targetRange2.Start + $markIdx
$endMark = $startMark + "@@MARK@@".Length
$rngMark = $d.Range($startMark, $endMark)
$rngMark.Text = ""

# Step C: re-split the now-fully-merged tail of the paragraph back into its
# original two runs (" for Windows operating systems" and the closing
# sentence) by toggling a formatting property on and back off across each
# desired run boundary -- this forces the host to re-partition the runs
# without leaving any net formatting change behind.
$target3 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text.StartsWith("You then must click the Download")) {
        $target3 = $cand
        break
    }
}
$targetRange3 = $target3.Range
$fullText3 = $targetRange3.Text

$marker4 = " for Windows operating systems"
$idx4 = $fullText3.IndexOf($marker4)
$start4 = $targetRange3.Start + $idx4
$end4 = $start4 + $marker4.Length
$rng4 = $d.Range($start4, $end4)
$rng4.Bold = 1
$rng4.Bold = 0

$marker5 = ". Once downloaded"
$idx5 = $fullText3.IndexOf($marker5)
$start5 = $targetRange3.Start + $idx5
$end5 = $targetRange3.End - 1
$rng5 = $d.Range($start5, $end5)
$rng5.Bold = 1
$rng5.Bold = 0

Write-Host "done"
